$d = $word.ActiveDocument
$d.Paragraphs(1).Range.Delete()
